$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(22, 23, 24, 25, 26, 27, 28, 29, 30, 31)
$aVals = @(-0.2223867177963257, 0.1144545078277593, 0.001384258270263117, -0.04791402816772423, 0.02978670597076406, 0.007494449615478475, 0.001274824142456048, 0.003990173339843651, -0.00699806213378908, 0.01654338836669934)
$bVals = @(0.8657857477664948, 0.9493652880191804, 0.9941152483224868, 0.9533933401107785, 0.9465968608856203, 0.9755845963954928, 0.8726469278335574, 0.9426385164260862, 0.904570586979389, 0.8796091377735137)
$cVals = @(0.5839412361383439, 0.8125896006822588, 0.7851106524467466, 0.6604413688182829, 0.6349874883890154, 0.702362984418869, 0.6991539821028709, 0.7098991274833678, 0.691669926047325, 0.693576380610466)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $aVals[$i]
    $ws.Cells.Item($r, 2).Value = $bVals[$i]
    $ws.Cells.Item($r, 3).Value = $cVals[$i]
}
